# Updated symbol list (crypto price feed snapshot refresh).
# For each coin row: Price (D), Volume(1h) % (E) and Hora (G) are refreshed to the
# new feed pull. Values are assigned with a leading apostrophe so Excel stores them
# as text (matching the workbook's existing text-typed cells) instead of auto-coercing
# numeric-looking strings into Number/Percentage cells; ClearFormats() then strips the
# resulting quote-prefix marker so the cell keeps its original default styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'265.36"
$ws.Range("E2").Value = "'1.56%"
$ws.Range("G2").Value = "'19"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").ClearFormats()
$ws.Range("G2").ClearFormats()

# Row 3
$ws.Range("D3").Value = "'26.77"
$ws.Range("E3").Value = "'-1.41%"
$ws.Range("G3").Value = "'19"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").ClearFormats()
$ws.Range("G3").ClearFormats()

# Row 4
$ws.Range("D4").Value = "'4.695"
$ws.Range("E4").Value = "'-0.30%"
$ws.Range("G4").Value = "'19"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").ClearFormats()
$ws.Range("G4").ClearFormats()

# Row 5
$ws.Range("D5").Value = "'0.06084"
$ws.Range("E5").Value = "'-1.75%"
$ws.Range("G5").Value = "'19"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").ClearFormats()
$ws.Range("G5").ClearFormats()

# Row 6
$ws.Range("D6").Value = "'6.744"
$ws.Range("E6").Value = "'0.43%"
$ws.Range("G6").Value = "'19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").ClearFormats()
$ws.Range("G6").ClearFormats()

# Row 7
$ws.Range("D7").Value = "'0.8504"
$ws.Range("E7").Value = "'0.03%"
$ws.Range("G7").Value = "'19"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").ClearFormats()
$ws.Range("G7").ClearFormats()

# Row 8
$ws.Range("D8").Value = "'0.9037"
$ws.Range("E8").Value = "'-1.19%"
$ws.Range("G8").Value = "'19"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").ClearFormats()
$ws.Range("G8").ClearFormats()

# Row 9
$ws.Range("D9").Value = "'0.1409"
$ws.Range("E9").Value = "'-0.19%"
$ws.Range("G9").Value = "'19"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").ClearFormats()
$ws.Range("G9").ClearFormats()

# Row 10
$ws.Range("D10").Value = "'0.04789"
$ws.Range("E10").Value = "'4.08%"
$ws.Range("G10").Value = "'19"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").ClearFormats()
$ws.Range("G10").ClearFormats()

# Row 11
$ws.Range("D11").Value = "'0.07101"
$ws.Range("E11").Value = "'0.23%"
$ws.Range("G11").Value = "'19"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").ClearFormats()
$ws.Range("G11").ClearFormats()

# Row 12
$ws.Range("D12").Value = "'0.03164"
$ws.Range("E12").Value = "'0.73%"
$ws.Range("G12").Value = "'19"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").ClearFormats()
$ws.Range("G12").ClearFormats()

# Row 13
$ws.Range("D13").Value = "'0.09018"
$ws.Range("E13").Value = "'-0.36%"
$ws.Range("G13").Value = "'19"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").ClearFormats()
$ws.Range("G13").ClearFormats()

# Row 14
$ws.Range("D14").Value = "'0.001535"
$ws.Range("G14").Value = "'19"
$ws.Range("D14").ClearFormats()
$ws.Range("G14").ClearFormats()

# Row 15
$ws.Range("D15").Value = "'0.0006070"
$ws.Range("E15").Value = "'-1.40%"
$ws.Range("G15").Value = "'19"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").ClearFormats()
$ws.Range("G15").ClearFormats()

# Row 16
$ws.Range("D16").Value = "'0.005997"
$ws.Range("E16").Value = "'-1.26%"
$ws.Range("G16").Value = "'19"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").ClearFormats()
$ws.Range("G16").ClearFormats()

# Row 17
$ws.Range("D17").Value = "'3.457"
$ws.Range("E17").Value = "'-0.08%"
$ws.Range("G17").Value = "'19"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").ClearFormats()
$ws.Range("G17").ClearFormats()

# Row 18
$ws.Range("D18").Value = "'3.168"
$ws.Range("E18").Value = "'0.11%"
$ws.Range("G18").Value = "'19"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").ClearFormats()
$ws.Range("G18").ClearFormats()

# Row 19
$ws.Range("D19").Value = "'2.278"
$ws.Range("E19").Value = "'3.82%"
$ws.Range("G19").Value = "'19"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").ClearFormats()
$ws.Range("G19").ClearFormats()

# Row 20
$ws.Range("E20").Value = "'-0.82%"
$ws.Range("G20").Value = "'19"
$ws.Range("E20").ClearFormats()
$ws.Range("G20").ClearFormats()

# Row 21
$ws.Range("D21").Value = "'0.1299"
$ws.Range("E21").Value = "'-0.84%"
$ws.Range("G21").Value = "'19"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").ClearFormats()
$ws.Range("G21").ClearFormats()

# Row 22
$ws.Range("D22").Value = "'4.088"
$ws.Range("E22").Value = "'-0.33%"
$ws.Range("G22").Value = "'19"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").ClearFormats()
$ws.Range("G22").ClearFormats()

# Row 23
$ws.Range("D23").Value = "'0.04235"
$ws.Range("E23").Value = "'0.06%"
$ws.Range("G23").Value = "'19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").ClearFormats()
$ws.Range("G23").ClearFormats()

# Row 24
$ws.Range("D24").Value = "'0.001184"
$ws.Range("E24").Value = "'-2.74%"
$ws.Range("G24").Value = "'19"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").ClearFormats()
$ws.Range("G24").ClearFormats()

# Row 25
$ws.Range("D25").Value = "'0.004130"
$ws.Range("E25").Value = "'8.63%"
$ws.Range("G25").Value = "'19"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").ClearFormats()
$ws.Range("G25").ClearFormats()

# Row 26
$ws.Range("G26").Value = "'19"
$ws.Range("G26").ClearFormats()

# Row 27
$ws.Range("E27").Value = "'5.00%"
$ws.Range("G27").Value = "'19"
$ws.Range("E27").ClearFormats()
$ws.Range("G27").ClearFormats()

# Row 28
$ws.Range("G28").Value = "'19"
$ws.Range("G28").ClearFormats()

# Row 29
$ws.Range("G29").Value = "'19"
$ws.Range("G29").ClearFormats()

# Row 30
$ws.Range("G30").Value = "'19"
$ws.Range("G30").ClearFormats()

# Row 31
$ws.Range("G31").Value = "'19"
$ws.Range("G31").ClearFormats()

# Row 32
$ws.Range("G32").Value = "'19"
$ws.Range("G32").ClearFormats()

# Row 33
$ws.Range("G33").Value = "'19"
$ws.Range("G33").ClearFormats()

# Row 34
$ws.Range("G34").Value = "'19"
$ws.Range("G34").ClearFormats()

# Row 35
$ws.Range("G35").Value = "'19"
$ws.Range("G35").ClearFormats()

# Row 36
$ws.Range("G36").Value = "'19"
$ws.Range("G36").ClearFormats()

# Row 37
$ws.Range("G37").Value = "'19"
$ws.Range("G37").ClearFormats()

# Row 38
$ws.Range("G38").Value = "'19"
$ws.Range("G38").ClearFormats()

# Row 39
$ws.Range("G39").Value = "'19"
$ws.Range("G39").ClearFormats()

# Row 40
$ws.Range("D40").Value = "'0.03914"
$ws.Range("E40").Value = "'-1.02%"
$ws.Range("G40").Value = "'19"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").ClearFormats()
$ws.Range("G40").ClearFormats()

# Row 41
$ws.Range("D41").Value = "'0.1115"
$ws.Range("E41").Value = "'0.19%"
$ws.Range("G41").Value = "'19"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").ClearFormats()
$ws.Range("G41").ClearFormats()

# Row 42
$ws.Range("D42").Value = "'0.004187"
$ws.Range("E42").Value = "'1.63%"
$ws.Range("G42").Value = "'19"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").ClearFormats()
$ws.Range("G42").ClearFormats()

# Row 43
$ws.Range("E43").Value = "'-3.33%"
$ws.Range("G43").Value = "'19"
$ws.Range("E43").ClearFormats()
$ws.Range("G43").ClearFormats()

# Row 44
$ws.Range("E44").Value = "'-16.84%"
$ws.Range("G44").Value = "'19"
$ws.Range("E44").ClearFormats()
$ws.Range("G44").ClearFormats()

# Row 45
$ws.Range("D45").Value = "'0.00005128"
$ws.Range("E45").Value = "'-0.55%"
$ws.Range("G45").Value = "'19"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").ClearFormats()
$ws.Range("G45").ClearFormats()

# Row 46
$ws.Range("E46").Value = "'0.00%"
$ws.Range("G46").Value = "'19"
$ws.Range("E46").ClearFormats()
$ws.Range("G46").ClearFormats()

# Row 47
$ws.Range("G47").Value = "'19"
$ws.Range("G47").ClearFormats()

# Row 48
$ws.Range("D48").Value = "'0.1585"
$ws.Range("E48").Value = "'-4.89%"
$ws.Range("G48").Value = "'19"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").ClearFormats()
$ws.Range("G48").ClearFormats()

# Row 49
$ws.Range("E49").Value = "'0.00%"
$ws.Range("G49").Value = "'19"
$ws.Range("E49").ClearFormats()
$ws.Range("G49").ClearFormats()

# Row 50
$ws.Range("E50").Value = "'0.00%"
$ws.Range("G50").Value = "'19"
$ws.Range("E50").ClearFormats()
$ws.Range("G50").ClearFormats()

# Row 51
$ws.Range("G51").Value = "'19"
$ws.Range("G51").ClearFormats()
